# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# to reflect refreshed values from the latest GitHub Actions scrape run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '19.921.29'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -8.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.410.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -8.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.002'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.00%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '273.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.84%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3702'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -5.33%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3072'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.11%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.20'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -8.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9973'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.55%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06570'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -8.62%  '

$ws.Range("E12").Value = '  -0.08%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.355'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.80%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.00'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -8.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.165'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.79%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.410.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -8.70%  '

$ws.Range("E17").Value = '  -8.45%  '

$ws.Range("E18").Value = '  -12.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '73.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -11.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.610'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.75%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.39'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.47%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.88'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.279'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.62%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '19.965.58'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -7.98%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.253'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '139.09'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.39%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -8.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.569.81'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.43%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '109.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.13%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.806'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -21.39%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.338'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.39%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8256'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -14.73%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07686'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.94%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.411'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.76%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05767'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.00%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.805'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.07%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.001'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.1934'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.68%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02045'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.20'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.49%  '

$ws.Range("E42").Value = '  -9.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.277'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -10.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5301'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.49%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.529'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.78%  '

$ws.Range("E46").Value = '  -6.50%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5113'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.67%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.805'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.36%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.74'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.69%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.042'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -10.20%  '

$ws.Range("E51").Value = '  -0.07%  '
